$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Runmode column (C) values: swap Y/N for several suites
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "Y"
$ws.Range("C6").Value = "N"

# Update the selected cell shown in the sheet view
$ws.Range("C5").Select()
